# Applies the weekly data refresh to the "Jengibre" (Ginger) price sheet.
# Two new weekly price records are inserted into the dataset:
#   - one inserted before the (old) row 37, becoming the new row 37
#   - one inserted before the (old) row 83 (which, after the first insertion,
#     sits at row 84), becoming the new row 84
# All subsequent rows are pushed down automatically by the row Insert,
# preserving their original values and cell formatting (e.g. the date
# number format on column D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow($Row, $Fecha, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg) {
    $ws.Cells.Item($Row, 1).Value = 8
    $ws.Cells.Item($Row, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($Row, 3).Value = "Coquimbo"
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = 4
    $ws.Cells.Item($Row, 6).Value = 100114007
    $ws.Cells.Item($Row, 7).Value = "Jengibre"
    $ws.Cells.Item($Row, 8).Value = "Sin especificar"
    $ws.Cells.Item($Row, 9).Value = "Primera"
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 14).Value = "`$/caja 13 kilos"
    $ws.Cells.Item($Row, 15).Value = "Perú"
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = 13
    $ws.Cells.Item($Row, 18).Value = "Hortaliza"
}

# Insert first new weekly record at row 37 (pushes old row37.. down by one)
$ws.Rows.Item(37).Insert()
Set-DataRow 37 44965 400 22500 23000 22750 1750

# Insert second new weekly record at row 84 (old row83 data, now at row84
# after the previous insert, gets pushed down to row85)
$ws.Rows.Item(84).Insert()
Set-DataRow 84 44964 300 23000 24000 23500 1808
